# Update "Elapsed Duration(Hrs)" values (column G) on several sheets to
# reflect a later recalculation time (commit: 6/18/2025, 11:24:57 AM).

$wb = $excel.ActiveWorkbook

$updates = @(
    @{ Sheet = "R1"; Cell = "G2"; Value = "3924:39:11" },
    @{ Sheet = "R1"; Cell = "G3"; Value = "64:11:49" },
    @{ Sheet = "R2"; Cell = "G2"; Value = "12106:02:52" },
    @{ Sheet = "R2"; Cell = "G3"; Value = "3235:46:21" },
    @{ Sheet = "R2"; Cell = "G4"; Value = "473:57:55" },
    @{ Sheet = "R4"; Cell = "G2"; Value = "2951:52:41" },
    @{ Sheet = "R4"; Cell = "G3"; Value = "179:04:56" },
    @{ Sheet = "R5"; Cell = "G2"; Value = "425:51:40" },
    @{ Sheet = "R6"; Cell = "G2"; Value = "66:23:58" }
)

foreach ($u in $updates) {
    $ws = $wb.Worksheets.Item($u.Sheet)
    $ws.Range($u.Cell).Value = $u.Value
}
